$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update 想去人数 (interest count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1317
$wsExhibit.Range("F3").Value = 2828

# Sheet "全部类型" (All types) - same two events appear here as well
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1317
$wsAll.Range("F4").Value = 2828
